$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 15 (nidio_spolis_month_2006_2023): update Unique Identifier text, File size, styling ---

$ws.Range("E15").Value = "37.1 GB"

$newId = "year-rinpersoon-baanrugid/ikvid" + " " + [char]10 + "year-rinpersoon (if only mainjob)"

$ws.Range("D15").Value = $newId
$ws.Range("D15").Characters(1, 31).Font.Bold = $true

# --- Row 16 (nidio_spolis_year_2006_2023): same Unique Identifier, clear File size ---

$ws.Range("D15").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws.Range("E16").ClearContents()

# --- Formatting: wrap the Unique Identifier column, vertically center the rest, taller rows ---

$ws.Range("D15:D16").WrapText = $true

$ws.Range("A15:C16").VerticalAlignment = -4108
$ws.Range("E15:F16").VerticalAlignment = -4108

$ws.Rows.Item(15).RowHeight = 30
$ws.Rows.Item(16).RowHeight = 30
